$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 15 from 45224 (2023-10-25)
# to 45233 (2023-11-03), keeping the existing date formatting.
$ws.Range("C2:C15").Value = 45233
